# "traded, fixed issues with the repeater"
# Row 3 (the GILD trade that was still open) now has its sell-side data
# filled in, and the position is recorded as unprofitable. A new trade
# (row 4) has been started with its starting principle.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Complete row 3: record the sale, its price-change %, and that the trade
# was not profitable.
$ws.Range("B3").Value = $false
$ws.Range("E3").Value = 75.5
$ws.Range("F3").Value = -1.2426383744658567
$ws.Range("G3").Value = $false

# Start row 4 for the next trade with the rolled-forward principle.
$ws.Range("C4").Value = 9814.58
